$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the 5th-choice value for project p3 (row 4) from "s9" to "s10"
$ws.Range("G4").Value = "s10"
